# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps on row 2 of the
# zh-cn and de-de sheets, as part of regenerating the handback status
# report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-21 17:09:04"
$wsZhCn.Range("H2").Value = "2016-03-21 17:09:25"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-21 17:09:08"
$wsDeDe.Range("H2").Value = "2016-03-21 17:09:30"
